# BOM.xlsx update: add SKBPC5004 three-phase bridge rectifier candidate row,
# turn the existing IGBT link into a real hyperlink, and flag the new row's
# "Labda var mı?" cell as unavailable (red) the same way the IGBT row is
# flagged as available (green).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row (row 3) link becomes a clickable hyperlink (adds the
# "Hyperlink" cell style / underlined themed font used by the workbook).
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.direnc.net/ixgh24n60c4d1-rohs-24a600v-to247ad-igbtdiode")

# New candidate component row (row 4).
$ws.Range("B4").Value = "SKBPC5004 Three Phase Bridge Rectifier"
$ws.Range("D4").Value = "Yok"
$ws.Range("D4").Interior.Color = 255
$ws.Range("F4").Value = "https://www.motorobit.com/skbpc5016-50a-1600v-trifaze-kopru-diyot-3-faz"
$ws.Range("G4").Value = "1600V"
$ws.Range("H4").Value = "50A"

# Leave selection where the author left it before saving.
$ws.Range("J7").Select() | Out-Null
